$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header capitalization: startDateTime -> startDatetime, endDateTime -> endDatetime
$ws.Range("C1").Value = "endDatetime"
$ws.Range("B1").Value = "startDatetime"

# Update the active selection to match the authored state
$ws.Range("G9").Select()
